# Updates to data_output/Excess Mortality/Countries/TOTAL_EUROPE_Total_...xlsx
# Row 7 (Female, Czechia)
# Row 14 (Female, Iceland)
# Row 36 (Male, Czechia)
# Row 38 (Male, Estonia)
# Row 43 (Male, Iceland)
# Row 65 (Total, Czechia)
# Row 67 (Total, Estonia)
# Row 72 (Total, Iceland)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: Female, Czechia ---
$ws.Range("D7").Value = 53445
$ws.Range("P7").Value = 8721.6
$ws.Range("U7").Value = "8721.6 (±750.6)"
$ws.Range("X7").Value = 160.9
$ws.Range("Z7").Value = "160.9(±13.8)"

# --- Row 14: Female, Iceland ---
$ws.Range("D14").Value = 914
$ws.Range("P14").Value = -5.2
$ws.Range("R14").Value = -0.6
$ws.Range("U14").Value = "-5.2 (±17.9)"
$ws.Range("V14").Value = "-0.6% (±1.9%)"
$ws.Range("X14").Value = -2.9
$ws.Range("Z14").Value = "-2.9(±10.1)"

# --- Row 36: Male, Czechia ---
$ws.Range("D36").Value = 56998
$ws.Range("P36").Value = 10605.2
$ws.Range("U36").Value = "10605.2 (±643.0)"
$ws.Range("X36").Value = 201.2
$ws.Range("Z36").Value = "201.2(±12.2)"

# --- Row 38: Male, Estonia ---
$ws.Range("D38").Value = 6328
$ws.Range("P38").Value = 329
$ws.Range("U38").Value = "329.0 (±68.8)"
$ws.Range("X38").Value = 52.3
$ws.Range("Y38").Value = 10.9
$ws.Range("Z38").Value = "52.3(±10.9)"

# --- Row 43: Male, Iceland ---
$ws.Range("D43").Value = 980
$ws.Range("P43").Value = 21.2
$ws.Range("R43").Value = 2.2
$ws.Range("U43").Value = "21.2 (±25.8)"
$ws.Range("V43").Value = "2.2% (±2.7%)"
$ws.Range("X43").Value = 11.3
$ws.Range("Z43").Value = "11.3(±13.8)"

# --- Row 65: Total, Czechia ---
$ws.Range("D65").Value = 110443
$ws.Range("P65").Value = 19326.8
$ws.Range("U65").Value = "19326.8 (±1350.5)"
$ws.Range("Y65").Value = 12.7
$ws.Range("Z65").Value = "180.7(±12.7)"

# --- Row 67: Total, Estonia ---
$ws.Range("D67").Value = 13414
$ws.Range("P67").Value = 771.2
$ws.Range("U67").Value = "771.2 (±116.7)"
$ws.Range("Y67").Value = 8.800000000000001
$ws.Range("Z67").Value = "58.0(±8.8)"

# --- Row 72: Total, Iceland ---
$ws.Range("D72").Value = 1894
$ws.Range("P72").Value = 16
$ws.Range("R72").Value = 0.9
$ws.Range("S72").Value = 2.2
$ws.Range("U72").Value = "16.0 (±40.9)"
$ws.Range("V72").Value = "0.9% (±2.2%)"
$ws.Range("X72").Value = 4.4
$ws.Range("Z72").Value = "4.4(±11.2)"
